$wb = $excel.ActiveWorkbook

# --- ALC (sheet1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4000
$ws.Range("J51").Value = 4000
$ws.Range("L51").Value = 4000
$ws.Range("N51").Value = -4968
$ws.Range("H92").Value = 799.35297
$ws.Range("I92").Value = 799.35297
$ws.Range("K92").Value = 799.35297
$ws.Range("M92").Value = 448.64703
$ws.Range("H113").Value = 2688.25
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 1254
$ws.Range("H125").Value = 3155
$ws.Range("I125").Value = 3298.3333
$ws.Range("J125").Value = 3011.6667
$ws.Range("K125").Value = 29684.9997
$ws.Range("L125").Value = 27105.0003
$ws.Range("M125").Value = -27224.9997
$ws.Range("N125").Value = -32025.0003
$ws.Range("H137").Value = 1269.2258
$ws.Range("J137").Value = 1396.1428
$ws.Range("L137").Value = 4188.428400000001
$ws.Range("N137").Value = -9288.428400000001
$ws.Range("H138").Value = 2057.5
$ws.Range("I138").Value = 782.35297
$ws.Range("J138").Value = 2318.6748
$ws.Range("K138").Value = 2347.05891
$ws.Range("L138").Value = 6956.024399999999
$ws.Range("M138").Value = 2792.94109
$ws.Range("N138").Value = -17236.0244

# --- ARM (sheet2) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2118.5
$ws.Range("I35").Value = 2118.5
$ws.Range("K35").Value = 2118.5
$ws.Range("M35").Value = -1712.5
$ws.Range("H37").Value = 14250
$ws.Range("I37").Value = 500
$ws.Range("K37").Value = 500
$ws.Range("M37").Value = -227
$ws.Range("H61").Value = 1549.2667
$ws.Range("I61").Value = 1441.4615
$ws.Range("K61").Value = 1441.4615
$ws.Range("M61").Value = -1229.4615
$ws.Range("H74").Value = 819.63336
$ws.Range("I74").Value = 848.5
$ws.Range("K74").Value = 848.5
$ws.Range("M74").Value = 25.5
$ws.Range("H77").Value = 819.63336
$ws.Range("I77").Value = 848.5
$ws.Range("K77").Value = 4242.5
$ws.Range("M77").Value = 125.5
$ws.Range("H136").Value = 1549.2667
$ws.Range("I136").Value = 1441.4615
$ws.Range("K136").Value = 4324.3845
$ws.Range("M136").Value = -1774.3845

# --- BSM (sheet3) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 10096.458
$ws.Range("I134").Value = 7500.9414
$ws.Range("J134").Value = 16399.857
$ws.Range("K134").Value = 22502.8242
$ws.Range("L134").Value = 49199.571
$ws.Range("M134").Value = -19967.8242
$ws.Range("N134").Value = -54269.571

# --- CRP (sheet4) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H127").Value = 49780
$ws.Range("J127").Value = 49780
$ws.Range("L127").Value = 49780
$ws.Range("N127").Value = -59700

# --- CUL (sheet5) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 932.1
$ws.Range("I6").Value = 81.25
$ws.Range("K6").Value = 243.75
$ws.Range("M6").Value = -130.75
$ws.Range("H10").Value = 40
$ws.Range("I10").Value = 40
$ws.Range("K10").Value = 120
$ws.Range("M10").Value = 19
$ws.Range("H17").Value = 1162.75
$ws.Range("J17").Value = 1162.75
$ws.Range("L17").Value = 3488.25
$ws.Range("N17").Value = -3826.25
$ws.Range("H132").Value = 2038.3846
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 13500
$ws.Range("M132").Value = -10970

# --- GSM (sheet6) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4120.091
$ws.Range("I80").Value = 1800
$ws.Range("J80").Value = 5445.857
$ws.Range("K80").Value = 1800
$ws.Range("L80").Value = 5445.857
$ws.Range("M80").Value = -802
$ws.Range("N80").Value = -7441.857
$ws.Range("H83").Value = 4120.091
$ws.Range("I83").Value = 1800
$ws.Range("J83").Value = 5445.857
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 27229.285
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -37213.285
$ws.Range("H126").Value = 2720.3635
$ws.Range("I126").Value = 1804.8
$ws.Range("J126").Value = 3483.3333
$ws.Range("K126").Value = 5414.4
$ws.Range("L126").Value = 10449.9999
$ws.Range("M126").Value = -2944.4
$ws.Range("N126").Value = -15389.9999
$ws.Range("H132").Value = 2118.9355
$ws.Range("I132").Value = 1727.24
$ws.Range("K132").Value = 5181.72
$ws.Range("M132").Value = -2651.72

# --- LTW (sheet7) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2157.3572
$ws.Range("I7").Value = 2024.75
$ws.Range("J7").Value = 2334.1667
$ws.Range("K7").Value = 2024.75
$ws.Range("L7").Value = 2334.1667
$ws.Range("M7").Value = -1912.75
$ws.Range("N7").Value = -2558.1667
$ws.Range("H16").Value = 1009.6087
$ws.Range("I16").Value = 952.4737
$ws.Range("J16").Value = 1281
$ws.Range("K16").Value = 952.4737
$ws.Range("L16").Value = 1281
$ws.Range("M16").Value = -782.4737
$ws.Range("N16").Value = -1621
$ws.Range("H22").Value = 1520
$ws.Range("J22").Value = 1577.7778
$ws.Range("L22").Value = 1577.7778
$ws.Range("N22").Value = -2167.7778
$ws.Range("H27").Value = 1520
$ws.Range("J27").Value = 1577.7778
$ws.Range("L27").Value = 1577.7778
$ws.Range("N27").Value = -1791.7778
$ws.Range("H40").Value = 2527
$ws.Range("I40").Value = 2338
$ws.Range("K40").Value = 2338
$ws.Range("M40").Value = -2202
$ws.Range("H46").Value = 983.3333
$ws.Range("J46").Value = 1225
$ws.Range("L46").Value = 1225
$ws.Range("N46").Value = -1601
$ws.Range("H68").Value = 1509
$ws.Range("I68").Value = 1259.3334
$ws.Range("J68").Value = 2133.1667
$ws.Range("K68").Value = 1259.3334
$ws.Range("L68").Value = 2133.1667
$ws.Range("M68").Value = -510.3334
$ws.Range("N68").Value = -3631.1667
$ws.Range("H71").Value = 1509
$ws.Range("I71").Value = 1259.3334
$ws.Range("J71").Value = 2133.1667
$ws.Range("K71").Value = 6296.666999999999
$ws.Range("L71").Value = 10665.8335
$ws.Range("M71").Value = -2552.666999999999
$ws.Range("N71").Value = -18153.8335
$ws.Range("H82").Value = 2299
$ws.Range("I82").Value = 2198.75
$ws.Range("K82").Value = 2198.75
$ws.Range("M82").Value = -1837.75
$ws.Range("H85").Value = 2299
$ws.Range("I85").Value = 2198.75
$ws.Range("K85").Value = 2198.75
$ws.Range("M85").Value = -950.75
$ws.Range("H126").Value = 2157.3572
$ws.Range("I126").Value = 2024.75
$ws.Range("J126").Value = 2334.1667
$ws.Range("K126").Value = 6074.25
$ws.Range("L126").Value = 7002.500100000001
$ws.Range("M126").Value = -3604.25
$ws.Range("N126").Value = -11942.5001

# --- WVR (sheet8) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1903.3429
$ws.Range("I132").Value = 1596.8
$ws.Range("J132").Value = 2669.7
$ws.Range("K132").Value = 4790.4
$ws.Range("L132").Value = 8009.099999999999
$ws.Range("M132").Value = -2260.4
$ws.Range("N132").Value = -13069.1
$ws.Range("H136").Value = 1900.9286
$ws.Range("I136").Value = 1713.125
$ws.Range("J136").Value = 2151.3333
$ws.Range("K136").Value = 5139.375
$ws.Range("L136").Value = 6453.999899999999
$ws.Range("M136").Value = -2589.375
$ws.Range("N136").Value = -11553.9999
$ws.Range("H141").Value = 36715
$ws.Range("J141").Value = 36715
$ws.Range("L141").Value = 36715
$ws.Range("N141").Value = -47075
